$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.82049999999999
$ws.Range("B9").Value = 8.643600000000005
$ws.Range("B18").Value = 4.556600000000003
$ws.Range("B20").Value = 5.728500000000001
$ws.Range("D21").Value = -7.434200000000005
